# Rewrite the test-case table (rows 2-22, columns A:C) to match the
# updated "Email / Password login" test-case spec. The edit reshuffles
# several test cases, rewords several steps/expected-results, and adds
# a brand-new 22nd row ("Login Attempt Limit Exceeded" follow-up step),
# which grows the used range from A1:C21 to A1:C22.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 21,3
$arr[0,0] = 'Verify Email and Password Fields Display'
$arr[0,1] = 'Navigate to the login page.'
$arr[0,2] = ''
$arr[1,0] = ''
$arr[1,1] = 'Verify that two text fields are present: one labeled ''Email'' and one labeled ''Password''.'
$arr[1,2] = 'Two text fields with the labels ''Email'' and ''Password'' are displayed.'
$arr[2,0] = 'Verify Password Field Masking'
$arr[2,1] = 'Navigate to the login page.'
$arr[2,2] = ''
$arr[3,0] = ''
$arr[3,1] = 'Enter text into the ''Password'' field.'
$arr[3,2] = 'The entered characters in the ''Password'' field are masked, not visible to the user.'
$arr[4,0] = 'Successful Login with Valid Credentials'
$arr[4,1] = 'Navigate to the login page.'
$arr[4,2] = ''
$arr[5,0] = ''
$arr[5,1] = 'Enter valid email and password into the respective fields.'
$arr[5,2] = ''
$arr[6,0] = ''
$arr[6,1] = 'Click the ''Login'' button.'
$arr[6,2] = 'The user is successfully redirected to the dashboard.'
$arr[7,0] = 'Successful Login with Mixed Case Credentials'
$arr[7,1] = 'Navigate to the login page.'
$arr[7,2] = ''
$arr[8,0] = ''
$arr[8,1] = 'Enter valid email and password with mixed case characters into the respective fields.'
$arr[8,2] = ''
$arr[9,0] = ''
$arr[9,1] = 'Click the ''Login'' button.'
$arr[9,2] = 'The user is successfully redirected to the dashboard.'
$arr[10,0] = 'Login Failure with Incorrect Email'
$arr[10,1] = 'Navigate to the login page.'
$arr[10,2] = ''
$arr[11,0] = ''
$arr[11,1] = 'Enter an incorrect email address into the ''Email'' field.'
$arr[11,2] = ''
$arr[12,0] = ''
$arr[12,1] = 'Enter a valid password into the ''Password'' field.'
$arr[12,2] = ''
$arr[13,0] = ''
$arr[13,1] = 'Click the ''Login'' button.'
$arr[13,2] = 'An error message is displayed informing the user that the email or password is incorrect.'
$arr[14,0] = 'Login Failure with Incorrect Password'
$arr[14,1] = 'Navigate to the login page.'
$arr[14,2] = ''
$arr[15,0] = ''
$arr[15,1] = 'Enter a valid email address into the ''Email'' field.'
$arr[15,2] = ''
$arr[16,0] = ''
$arr[16,1] = 'Enter an incorrect password into the ''Password'' field.'
$arr[16,2] = ''
$arr[17,0] = ''
$arr[17,1] = 'Click the ''Login'' button.'
$arr[17,2] = 'An error message is displayed informing the user that the email or password is incorrect.'
$arr[18,0] = 'Login Attempt Limit Exceeded'
$arr[18,1] = 'Navigate to the login page.'
$arr[18,2] = ''
$arr[19,0] = ''
$arr[19,1] = 'Attempt to login with incorrect credentials 5 times within an hour.'
$arr[19,2] = ''
$arr[20,0] = ''
$arr[20,1] = 'Attempt to login again with incorrect credentials.'
$arr[20,2] = 'An error message is displayed informing the user that they have exceeded the login attempt limit. The message should also provide an estimated time for the restriction to be lifted.'

$ws.Range("A2:C22").Value = $arr

